# "Ajustando as bases para analise"
# The ProdutoID column (A2:A61) previously held raw GUID values; they are
# replaced here with short, readable product codes P001..P060 (row 2 -> P001,
# row 3 -> P002, ... row 61 -> P060), keeping every other column untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow  = 61

for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
    $n = $row - 1
    $code = "P{0:D3}" -f $n
    $ws.Cells.Item($row, 1).Value = $code
}

# Reflect the author's on-screen state after the edit: the new ProdutoID
# column is selected and the view has scrolled down a bit.
[void]$ws.Range("A2:A61").Select()
$excel.ActiveWindow.ScrollRow = 7
